$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table (rows 31-40) with the 10 columns extracted from the
# 'description' field: ambientes, baños, dormitorios (str) and
# balcon, cochera, jardin, lavadero, parrilla, patio, pileta (bool).
# Values are written in the exact order the strings were first
# authored so the shared-string table comes out in the same order.
$ws.Range("B31").Value = "agregados"
$ws.Range("D31").Value = "ambientes"
$ws.Range("I31").Value = "10292 no nulos - donde rooms es nulo"
$ws.Range("D32").Value = "baños"
$ws.Range("D33").Value = "dormitorios"
$ws.Range("D34").Value = "balcon"
$ws.Range("F34").Value = "bool"
$ws.Range("F31").Value = "str - pero es un número"
$ws.Range("D35").Value = "cochera"
$ws.Range("I35").Value = "41977 True"
$ws.Range("I34").Value = "37149 True"
$ws.Range("I33").Value = "34366 no nulos"
$ws.Range("I32").Value = "9965 no nulos"
$ws.Range("D36").Value = "jardin"
$ws.Range("I36").Value = "15334 True"
$ws.Range("D37").Value = "lavadero"
$ws.Range("I37").Value = "32707 True"
$ws.Range("D38").Value = "parrilla"
$ws.Range("I38").Value = "29778 True"
$ws.Range("D39").Value = "patio"
$ws.Range("I39").Value = "17269 True"
$ws.Range("D40").Value = "pileta"
$ws.Range("I40").Value = "28821 True"
$ws.Range("F32").Value = "str - pero es un número"
$ws.Range("F33").Value = "str - pero es un número"
$ws.Range("F35").Value = "bool"
$ws.Range("F36").Value = "bool"
$ws.Range("F37").Value = "bool"
$ws.Range("F38").Value = "bool"
$ws.Range("F39").Value = "bool"
$ws.Range("F40").Value = "bool"

# Restore the view state (selection) left by the author.
$ws.Range("D10").Select()
